$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.268.87"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "2.421.61"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.71%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +1.55%  "
$ws.Range("D9").Value = "2.419.68"
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.110"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.88%  "
$ws.Range("E11").Value = "  -1.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000176"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").Value = "2.859.61"
$ws.Range("E16").Value = "  +1.88%  "
$ws.Range("D17").Value = "62.136.04"
$ws.Range("E17").Value = "  +1.16%  "
$ws.Range("D18").Value = "2.419.56"
$ws.Range("E18").Value = "  +1.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.72%  "
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.92%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.45%  "
$ws.Range("E25").Value = "  -2.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "579.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.66%  "
$ws.Range("D28").Value = "2.540.73"
$ws.Range("E28").Value = "  +1.85%  "
$ws.Range("E29").Value = "  +0.45%  "
$ws.Range("D30").Value = "0.0₃0948"
$ws.Range("E30").Value = "  +5.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.150"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("E34").Value = "  +2.20%  "
$ws.Range("E35").Value = "  +2.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.72"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.14%  "
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.79"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.384"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "151.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.02%  "
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("E42").Value = "  -4.40%  "
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "150.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.59%  "
$ws.Range("E47").Value = "  +2.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.20%  "
$ws.Range("E49").Value = "  +2.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0922"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.00%  "
$ws.Range("E51").Value = "  +2.16%  "
